$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns B, C, E: plain text values (non-numeric strings, safe to assign directly)
$ws.Range('B19').Value = 'Chainlink'
$ws.Range('B20').Value = 'ShibaInu'
$ws.Range('B44').Value = 'MXToken'
$ws.Range('B45').Value = 'mCoin'
$ws.Range('B50').Value = 'Algorand'
$ws.Range('B51').Value = 'BabyDogeCoin'
$ws.Range('C19').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('C20').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('C44').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('C45').Value = 'https://coinranking.com/coin/fzVgyjBcRc9+mcoin-mcoin'
$ws.Range('C50').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('C51').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('E2').Value = '  +0.62%  '
$ws.Range('E3').Value = '  +0.66%  '
$ws.Range('E4').Value = '  -0.74%  '
$ws.Range('E5').Value = '  -0.11%  '
$ws.Range('E6').Value = '  +0.03%  '
$ws.Range('E7').Value = '  -0.86%  '
$ws.Range('E8').Value = '  +2.04%  '
$ws.Range('E9').Value = '  +0.97%  '
$ws.Range('E10').Value = '  +0.22%  '
$ws.Range('E11').Value = '  -2.98%  '
$ws.Range('E12').Value = '  +0.33%  '
$ws.Range('E13').Value = '  +0.13%  '
$ws.Range('E14').Value = '  +0.46%  '
$ws.Range('E15').Value = '  +1.64%  '
$ws.Range('E16').Value = '  +1.75%  '
$ws.Range('E17').Value = '  +0.34%  '
$ws.Range('E18').Value = '  +1.05%  '
$ws.Range('E19').Value = '  +1.13%  '
$ws.Range('E20').Value = '  -0.21%  '
$ws.Range('E21').Value = '  -0.39%  '
$ws.Range('E22').Value = '  +7.76%  '
$ws.Range('E23').Value = '  +2.20%  '
$ws.Range('E24').Value = '  +2.55%  '
$ws.Range('E25').Value = '  +0.22%  '
$ws.Range('E26').Value = '  +0.02%  '
$ws.Range('E27').Value = '  +0.91%  '
$ws.Range('E28').Value = '  +0.00%  '
$ws.Range('E29').Value = '  -0.74%  '
$ws.Range('E30').Value = '  +0.07%  '
$ws.Range('E31').Value = '  +0.11%  '
$ws.Range('E32').Value = '  +0.14%  '
$ws.Range('E33').Value = '  +0.49%  '
$ws.Range('E34').Value = '  +0.44%  '
$ws.Range('E35').Value = '  +0.87%  '
$ws.Range('E36').Value = '  -2.08%  '
$ws.Range('E37').Value = '  +2.51%  '
$ws.Range('E38').Value = '  +1.75%  '
$ws.Range('E39').Value = '  +0.14%  '
$ws.Range('E40').Value = '  -1.50%  '
$ws.Range('E41').Value = '  -0.20%  '
$ws.Range('E42').Value = '  -0.30%  '
$ws.Range('E43').Value = '  -0.57%  '
$ws.Range('E44').Value = '  +2.49%  '
$ws.Range('E45').Value = '  -0.31%  '
$ws.Range('E46').Value = '  -0.81%  '
$ws.Range('E47').Value = '  +0.25%  '
$ws.Range('E48').Value = '  +3.73%  '
$ws.Range('E49').Value = '  +2.27%  '
$ws.Range('E50').Value = '  +1.34%  '
$ws.Range('E51').Value = '  -0.25%  '

# Column D: force text storage to match original inlineStr (avoid numeric auto-conversion)
$dCells = @('D2', 'D3', 'D4', 'D5', 'D6', 'D7', 'D8', 'D10', 'D11', 'D12', 'D13', 'D14', 'D15', 'D16', 'D17', 'D18', 'D19', 'D20', 'D21', 'D22', 'D23', 'D25', 'D26', 'D27', 'D29', 'D31', 'D32', 'D33', 'D34', 'D35', 'D37', 'D38', 'D39', 'D40', 'D41', 'D43', 'D44', 'D45', 'D46', 'D47', 'D48', 'D49', 'D50', 'D51')
$dValues = @{
    'D2' = '27.845.84'
    'D3' = '1.641.09'
    'D4' = '0.997'
    'D5' = '212.03'
    'D6' = '0.524'
    'D7' = '0.996'
    'D8' = '23.37'
    'D10' = '0.0613'
    'D11' = '0.0864'
    'D12' = '1.868.30'
    'D13' = '1.636.01'
    'D14' = '4.06'
    'D15' = '0.564'
    'D16' = '65.50'
    'D17' = '27.788.88'
    'D18' = '232.46'
    'D19' = '7.65'
    'D20' = '0.0₃0720'
    'D21' = '0.999'
    'D22' = '10.72'
    'D23' = '4.40'
    'D25' = '150.12'
    'D26' = '6.91'
    'D27' = '15.70'
    'D29' = '0.997'
    'D31' = '0.0483'
    'D32' = '3.30'
    'D33' = '1.470.44'
    'D34' = '3.10'
    'D35' = '1.56'
    'D37' = '0.934'
    'D38' = '0.885'
    'D39' = '0.0168'
    'D40' = '0.558'
    'D41' = '69.17'
    'D43' = '0.998'
    'D44' = '2.28'
    'D45' = '2.45'
    'D46' = '5.37'
    'D47' = '1.779.14'
    'D48' = '1.76'
    'D49' = '87.84'
    'D50' = '0.100'
    'D51' = '0.0₇0997'
}
foreach ($addr in $dCells) {
    $rng = $ws.Range($addr)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $dValues[$addr]
    $rng.Style = $origStyle
}

Write-Output "Applied all changes"